$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample")

$ws.Range("B2").Value = "extraction protocol"
$ws.Range("C2").Value = "EFO"
$ws.Range("D2").Value = "http://purl.obolibrary.org/obo/EFO_0000490"
$ws.Range("E2").Value = "RNA_Extraction.txt"
$ws.Range("F2").Value = "Zea mays"
$ws.Range("G2").Value = "NCBITaxon"
$ws.Range("H2").Value = "http://purl.obolibrary.org/obo/NCBITaxon_4577"
$ws.Range("I2").Value = "leaf"
$ws.Range("J2").Value = "PO"
$ws.Range("K2").Value = "http://purl.obolibrary.org/obo/PO_0025034"
$ws.Range("L2").Value = "total RNA"
$ws.Range("M2").Value = "EFO"
$ws.Range("N2").Value = "http://purl.obolibrary.org/obo/EFO_0004964"
